$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Insert the new "welcome" paragraph text in front of the paragraph
# that is currently the empty, unformatted paragraph (4th paragraph). Typing
# into that location keeps the inserted text free of the NoSpacing / sz30 /
# shd formatting that the header (Andy Huang / Contact / Email) paragraphs
# carry.
# ---------------------------------------------------------------------------
$welcomeText = "I want to extend a warm welcome to {receiver_name} and invite them to attend {event_name}. It will be hosted at {venue_name} at the lovely time of {event_time}. I hope to see you there and the dress code is {dress_code}. Once again I look forward to seeing you there and to have a wonderful evening.`r"
$p4 = $d.Paragraphs.Item(4)
$insPt = $d.Range($p4.Range.Start, $p4.Range.Start)
$insPt.InsertBefore($welcomeText)

# ---------------------------------------------------------------------------
# Step 2: Remove the old header paragraphs -- "Andy Huang", "Contact: ...",
# "Email: ..." -- which carried the NoSpacing style / sz30 / shading that we
# don't want in the new document.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(3).Range.Delete()
$d.Paragraphs.Item(2).Range.Delete()
$d.Paragraphs.Item(1).Range.Delete()

# ---------------------------------------------------------------------------
# Step 3: Drop the "Sincerely," and "Andy Huang" sign-off paragraphs, along
# with one of the two blank paragraphs that trailed them. We select from the
# start of "Sincerely," through to the start of the very last (document
# terminating) empty paragraph and delete that whole range -- this merges
# everything down into that final paragraph mark (which can't itself be
# deleted) leaving a single empty paragraph behind.
# ---------------------------------------------------------------------------
$sincerelyPara = $d.Paragraphs.Item(6)
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Range($sincerelyPara.Range.Start, $lastPara.Range.Start).Delete()

# ---------------------------------------------------------------------------
# Step 4: Remove the long "I am super eager..." paragraph entirely.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(4).Range.Delete()

# ---------------------------------------------------------------------------
# Step 5: Turn the old "I am interested..." paragraph into "Kind regards,"
# ---------------------------------------------------------------------------
$kindRegardsPara = $d.Paragraphs.Item(3)
$kindRegardsRange = $kindRegardsPara.Range
[void]$kindRegardsRange.MoveEnd(1, -1)
$kindRegardsRange.Text = "Kind regards,"

# ---------------------------------------------------------------------------
# Step 6: Merge the old "I believe I am a qualified candidate..." paragraph
# together with the final (document terminating) empty paragraph, the same
# way as step 3, then type "{senders_name}" into what is now the last
# paragraph so it ends up unformatted and is the very last paragraph in the
# body (immediately followed by the sectPr).
# ---------------------------------------------------------------------------
$believeParaStart = $d.Paragraphs.Item(4)
$finalEmptyPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Range($believeParaStart.Range.Start, $finalEmptyPara.Range.Start).Delete()

$sendersPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$sendersIns = $d.Range($sendersPara.Range.Start, $sendersPara.Range.Start)
$sendersIns.InsertBefore("{senders_name}")
